$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A, D, E, F, G, H, Q, R swap their values between row 16 and row 17.
# Column B gets new distinct values (not a simple swap).

$swapCols = @("A", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $swapCols) {
    $rng16 = $ws.Range("$col" + "16")
    $rng17 = $ws.Range("$col" + "17")
    $v16 = $rng16.Value2
    $v17 = $rng17.Value2
    $rng16.Value = $v17
    $rng17.Value = $v16
}

# Column B: distinct new values
$ws.Range("B16").Value = 98980
$ws.Range("B17").Value = 81711
